$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Order matters: it controls the order new shared strings are appended in,
# which must match the target sharedStrings.xml table layout.

# 1) NumPoliza on row 3 becomes text "04104013002" (quote-prefixed so the
#    leading zero is preserved) -> claims shared-string index 14.
$ws.Range("E3").Value = "'04104013002"

# 2) Ambiente column (A) now points at the "i-" prefixed preprod host
#    -> claims shared-string index 15.
$ws.Range("A2").Value = "i-preproducciongestion.segurossura.com.ar"
$ws.Range("A3").Value = "i-preproducciongestion.segurossura.com.ar"

# 3) URL column (B, hyperlinked) now points at the "i-" prefixed preprod
#    PolicyCenter URL -> claims shared-string index 16.
$ws.Range("B2").Value = "https://i-preproducciongestion.segurossura.com.ar/pc/PolicyCenter.do"
$ws.Range("B3").Value = "https://i-preproducciongestion.segurossura.com.ar/pc/PolicyCenter.do"

# 4) NumPoliza on row 2 becomes a plain numeric literal.
$ws.Range("E2").Value = 12112001742

# Move the active selection to B4, matching the saved sheet view state.
$ws.Range("B4").Select()
